$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.647.81'
$ws.Range("E2").Value = '  +2.48%  '

$ws.Range("D3").Value = '1.677.02'
$ws.Range("E3").Value = '  +2.75%  '

$ws.Range("E4").Value = '  -0.22%  '

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '219.18'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  +2.12%  '

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '0.529'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  +1.95%  '

$ws.Range("E7").Value = '  -0.19%  '

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '29.14'
$cell.Style = "Normal"
$ws.Range("E8").Value = '  +1.70%  '

$ws.Range("E9").Value = '  +2.25%  '

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '0.0643'
$cell.Style = "Normal"
$ws.Range("E10").Value = '  +5.71%  '

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '0.0903'
$cell.Style = "Normal"
$ws.Range("E11").Value = '  -0.07%  '

$ws.Range("D12").Value = '1.917.88'
$ws.Range("E12").Value = '  +2.73%  '

$ws.Range("D13").Value = '1.669.23'
$ws.Range("E13").Value = '  +2.32%  '

$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '10.14'
$cell.Style = "Normal"
$ws.Range("E14").Value = '  +9.30%  '

$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '0.603'
$cell.Style = "Normal"
$ws.Range("E15").Value = '  +7.36%  '

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '4.01'
$cell.Style = "Normal"
$ws.Range("E16").Value = '  +4.57%  '

$ws.Range("D17").Value = '30.652.71'
$ws.Range("E17").Value = '  +2.44%  '

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '65.99'
$cell.Style = "Normal"
$ws.Range("E18").Value = '  +3.01%  '

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '243.02'
$cell.Style = "Normal"
$ws.Range("E19").Value = '  +0.82%  '

$ws.Range("D20").Value = '0.0₃0719'
$ws.Range("E20").Value = '  +2.61%  '

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$cell.Style = "Normal"

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '4.24'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  +2.69%  '

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '9.97'
$cell.Style = "Normal"
$ws.Range("E23").Value = '  +2.32%  '

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '2.16'
$cell.Style = "Normal"
$ws.Range("E24").Value = '  -0.70%  '

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '159.34'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  +0.66%  '

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '15.82'
$cell.Style = "Normal"
$ws.Range("E26").Value = '  +2.49%  '

$ws.Range("E27").Value = '  +2.27%  '

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '6.69'
$cell.Style = "Normal"
$ws.Range("E28").Value = '  +1.90%  '

$ws.Range("E29").Value = '  -0.25%  '

$ws.Range("E30").Value = '  +0.86%  '

$ws.Range("E31").Value = '  +3.94%  '

$ws.Range("E32").Value = '  +2.56%  '

$ws.Range("E33").Value = '  +3.90%  '

$ws.Range("D34").Value = '1.506.81'
$ws.Range("E34").Value = '  +5.88%  '

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '1.78'
$cell.Style = "Normal"
$ws.Range("E35").Value = '  +6.77%  '

$ws.Range("B36").Value = 'Aave'
$ws.Range("C36").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '83.21'
$cell.Style = "Normal"
$ws.Range("E36").Value = '  +10.36%  '

$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '1.02'
$cell.Style = "Normal"
$ws.Range("E37").Value = '  -0.50%  '

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '0.598'
$cell.Style = "Normal"
$ws.Range("E38").Value = '  +8.04%  '

$ws.Range("E39").Value = '  +4.54%  '

$ws.Range("E40").Value = '  -3.04%  '

$ws.Range("E41").Value = '  -0.02%  '

$ws.Range("E42").Value = '  +1.80%  '

$ws.Range("E43").Value = '  +1.48%  '

$ws.Range("E44").Value = '  +0.13%  '

$ws.Range("E45").Value = '  +1.14%  '

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.Style = "Normal"
$ws.Range("E46").Value = '  -0.12%  '

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '5.55'
$cell.Style = "Normal"
$ws.Range("E47").Value = '  +4.27%  '

$ws.Range("D48").Value = '1.808.87'
$ws.Range("E48").Value = '  +1.91%  '

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '49.76'
$cell.Style = "Normal"
$ws.Range("E49").Value = '  -1.84%  '

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '92.89'
$cell.Style = "Normal"
$ws.Range("E50").Value = '  +2.61%  '

$ws.Range("E51").Value = '  +4.10%  '
